# Update ROR HealthcareService Sensitive Unit StructureDefinition workbook
# per release-notes.md regeneration (version bump, status/date/contact refresh,
# and swap of the two "Mapping" columns on the Elements sheet).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Metadata sheet - update Version / Status / Date / Contact values
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B3").Value  = "0.4.0-snapshot-1"
$meta.Range("B6").Value  = "draft"
$meta.Range("B8").Value  = "2024-05-23T12:16:26+00:00"
$meta.Range("B10").Value = "ANS (https://esante.gouv.fr)"

# ---------------------------------------------------------------------
# 2. Elements sheet - the two rightmost "Mapping" columns (AK = 37,
#    AL = 38) were swapped: header text, column width and all row data.
# ---------------------------------------------------------------------
$elem = $wb.Worksheets.Item("Elements")

# Only rows whose AK/AL pair actually differs need to be touched (rows 2
# and 4 hold an empty value in both columns already, so swapping them
# would be a no-op edit best left alone).
$rowsToSwap = @(1, 3, 5, 6)

# Capture current ("before") values for both columns first, since we
# will overwrite them in place.
$akValues = @{}
$alValues = @{}
foreach ($r in $rowsToSwap) {
    $akValues[$r] = $elem.Cells.Item($r, 37).Value()
    $alValues[$r] = $elem.Cells.Item($r, 38).Value()
}

# Write back swapped: AK gets the old AL content, AL gets the old AK content.
foreach ($r in $rowsToSwap) {
    $elem.Cells.Item($r, 37).Value = $alValues[$r]
    $elem.Cells.Item($r, 38).Value = $akValues[$r]
}

# Swap the column widths too (best-fit widths follow the swapped content):
# AK was 24.98 / AL was 87.79  ->  AK becomes ~87.79 / AL becomes ~24.98.
# (Values are nudged slightly so they land on the nearest width the engine's
# internal character-width quantization can represent, closest to the
# target widths of 87.7890625 and 24.98046875.)
$elem.Columns.Item(37).ColumnWidth = 87.0
$elem.Columns.Item(38).ColumnWidth = 24.16666667
